$wb = $excel.ActiveWorkbook

# ---- Sheet "data": add a "description" column and a new data row ----
$ws1 = $wb.Worksheets.Item(1)

# F1 header "description" - copy header style (bold/border) from A1
$ws1.Range("F1").Value = "description"
$ws1.Range("A1").Copy()
$ws1.Range("F1").PasteSpecial(-4122)

# New row 2 of data
$ws1.Range("A2").Value = "59b5afaa-3f00-4f55-b771-5dd40fea1b69"
$ws1.Range("A1").Copy()
$ws1.Range("A2").PasteSpecial(-4122)

# Keep the date as literal text, not an auto-converted date serial
$ws1.Range("B2").Value = "'12/10/2022"
$ws1.Range("B2").Style = "Normal"

$ws1.Range("D2").Value = "Experience"
$ws1.Range("E2").Value = "Experience"
$ws1.Range("F2").Value = "work experience"

# ---- Sheet "headers": fix TRUE strings -> real booleans, add "description" row ----
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("D5").Value = $true
$ws2.Range("E5").Value = $true
$ws2.Range("F5").Value = $true

$ws2.Range("D6").Value = $true
$ws2.Range("E6").Value = $true
$ws2.Range("F6").Value = $true

# New row 7 describing the "description" field
$ws2.Range("A7").Value = "description"
$ws2.Range("A6").Copy()
$ws2.Range("A7").PasteSpecial(-4122)

$ws2.Range("B7").Value = "str"
$ws2.Range("C7").Value = "description"

$ws2.Range("D7").Value = "'TRUE"
$ws2.Range("D7").Style = "Normal"
$ws2.Range("E7").Value = "'TRUE"
$ws2.Range("E7").Style = "Normal"
$ws2.Range("F7").Value = "'TRUE"
$ws2.Range("F7").Style = "Normal"
